$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fill in the previously-empty "xt30 f" row (row 3)
$ws.Range("B3").Value = "xt30 f"
$ws.Range("C3").Value = 5

# Add new row 4: "xt30 m"
$ws.Range("B4").Value = "xt30 m"
$ws.Range("C4").Value = 5

# Add new row 5: "Conn_01x06_Pin"
$ws.Range("B5").Value = "Conn_01x06_Pin"
$ws.Range("C5").Value = 15

# Match the saved selection/active cell from the diff
$ws.Range("C5").Select()
